$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 7 ("Solicita-se acesso ao Processo nº 0032185265/2024.")
#    This shifts the old row 8 (Inteligencia Artificial...) up to row 7,
#    old row 9 (Casei recentemente...) up to row 8, etc.
$ws.Rows(7).Delete()

# 2) Insert a brand-new row at position 9 (right after "Casei recentemente...",
#    which is now row 8) to hold a new survey response.
$ws.Rows(9).Insert()

$ws.Range("A9").Value = 62
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 'Bom dia,   Eu ja tenho a viabilidade aprovada. porém a junta me pede para eu fazer uma nova viabilidade alegando o nome da empresa!  neste caso o número protocolo integrado vai mudar e a taxa foi paga como faço para aproveitamento a viabilidade aprovada e só atualizar o nome e as exigência da junta comercial. Pois pelo que eu entendo para eu fazer uma nova viabilidade tenho que cancelar a que está vigente, fazer uma nova com certeza terá outro número e como vou fazer com a taxa paga.  motivo pendência nota explicativa 1. corrigir nire: 7893214568-7 89- prezado senhor usuário, orientamos fazer uma nova viabilidade de nome empresarial e retirar o ( 01 ) constante após a natureza jurídica do nome empresarial  CO S DE E Ltda  Preciso de orientação em referencia a taxa que foi pago com o protocolo DFP4568523652, COMO MANTER ENTE NUMERO'
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = $false

# 3) Delete the BIOCASA row (still at row 15 after the delete+insert above,
#    since those two operations cancel out in row count up to this point).
$ws.Rows(15).Delete()
